# Burndown.xlsx update: "Pause Menu, Character Objects, Play Pause Step Fast
# Forward and Stop Game"
#
# Sprint 2 (rows 16-31) burndown: the "Actual" column (D) had stalled at
# row 27 (the team hit 0 remaining). Fill in the trailing zero actuals for
# the rest of the sprint (rows 28-31) so the burndown line continues flat
# at 0 through to the end of the sprint.
#
# Sprint 3 (rows 33-48) burndown: fill in "Actual" values for the days that
# have now passed (rows 38-41), and flip the highlight color on the three
# backlog items that are now complete -- "Create Character objects in
# game" (row 34... already had D34 set, only style elsewhere), "Pause,
# Stop, and Step Simulation" (row 36) and "Pause Menu" (row 39) -- from the
# "not started" orange fill (style 9) to the "done" fill (style 2), the
# same style already used on finished rows like G35/G40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 2 burndown: trailing "Actual" values (D28:D31) -------------
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0

# --- Sprint 3 burndown: newly recorded "Actual" values (D38:D41) -------
$ws.Range("D38").Value = 5
$ws.Range("D39").Value = 5
$ws.Range("D40").Value = 4
$ws.Range("D41").Value = 2

# --- Sprint 3 backlog: mark completed items (copy the "done" fill from -
# --- a cell that already uses it, e.g. G35, onto the newly finished -----
# --- rows) --------------------------------------------------------------
$ws.Range("G35").Copy()
$ws.Range("G34").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G36").PasteSpecial(-4122)
$ws.Range("G39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the saved view/selection to match where the user ended up --
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("E48").Select()

# --- Keep the embedded burndown charts in sync with the new data -------
for ($i = 1; $i -le $ws.ChartObjects().Count; $i++) {
    $co = $ws.ChartObjects($i)
    $chart = $co.Chart
    for ($j = 1; $j -le $chart.SeriesCollection().Count; $j++) {
        $ser = $chart.SeriesCollection($j)
        $ser.Formula = $ser.Formula
    }
    $chart.Refresh()
}

$wb.Save()
